# Updates cryptos list values (price + 1h volume change) per source diff.
# Cells are forced to Text format before/through the write so that numeric-
# looking strings (e.g. "292.12", "0.0710") are not auto-converted to numbers
# and lose formatting (trailing zeros, grouping dots, etc.); the style index
# is then reset to "Normal" so no stray cell styles are introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '39.785.05'
Set-TextValue 'E2' '  -0.06%  '
Set-TextValue 'D3' '2.214.74'
Set-TextValue 'E3' '  -0.25%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '292.12'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'D6' '87.11'
Set-TextValue 'E6' '  +1.02%  '
Set-TextValue 'D7' '0.511'
Set-TextValue 'E7' '  -0.73%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'D9' '0.468'
Set-TextValue 'E9' '  -0.98%  '
Set-TextValue 'D10' '30.62'
Set-TextValue 'E10' '  -0.73%  '
Set-TextValue 'D11' '0.0778'
Set-TextValue 'E11' '  -0.94%  '
Set-TextValue 'D12' '49.98'
Set-TextValue 'E12' '  +6.01%  '
Set-TextValue 'D13' '0.111'
Set-TextValue 'E13' '  +2.78%  '
Set-TextValue 'D14' '6.45'
Set-TextValue 'E14' '  +1.51%  '
Set-TextValue 'D15' '2.563.19'
Set-TextValue 'E15' '  -0.01%  '
Set-TextValue 'B16' 'WrappedEther'
Set-TextValue 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '2.259.91'
Set-TextValue 'E16' '  +1.81%  '
Set-TextValue 'B17' 'Chainlink'
Set-TextValue 'C17' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D17' '13.72'
Set-TextValue 'E17' '  -2.76%  '
Set-TextValue 'D18' '0.731'
Set-TextValue 'E18' '  +0.17%  '
Set-TextValue 'D19' '39.780.10'
Set-TextValue 'E19' '  +0.04%  '
Set-TextValue 'D20' '0.0₃0884'
Set-TextValue 'E20' '  +0.13%  '
Set-TextValue 'D21' '11.22'
Set-TextValue 'E21' '  +1.23%  '
Set-TextValue 'D22' '5.76'
Set-TextValue 'E22' '  -0.99%  '
Set-TextValue 'D23' '65.57'
Set-TextValue 'E23' '  -0.04%  '
Set-TextValue 'D24' '237.11'
Set-TextValue 'E24' '  +0.45%  '
Set-TextValue 'D25' '0.999'
Set-TextValue 'E25' '  -0.22%  '
Set-TextValue 'D26' '2.45'
Set-TextValue 'E26' '  -0.62%  '
Set-TextValue 'D27' '1.84'
Set-TextValue 'E27' '  -0.30%  '
Set-TextValue 'D28' '23.43'
Set-TextValue 'E28' '  +2.99%  '
Set-TextValue 'E29' '  -2.68%  '
Set-TextValue 'D30' '9.24'
Set-TextValue 'E30' '  -0.23%  '
Set-TextValue 'D31' '156.76'
Set-TextValue 'E31' '  +3.47%  '
Set-TextValue 'D32' '31.89'
Set-TextValue 'E32' '  -3.18%  '
Set-TextValue 'D33' '0.998'
Set-TextValue 'E33' '  +0.03%  '
Set-TextValue 'D34' '4.95'
Set-TextValue 'E34' '  -0.10%  '
Set-TextValue 'D35' '0.0710'
Set-TextValue 'E35' '  -1.29%  '
Set-TextValue 'E36' '  +3.93%  '
Set-TextValue 'D37' '2.33'
Set-TextValue 'E37' '  -1.81%  '
Set-TextValue 'E38' '  -0.44%  '
Set-TextValue 'D39' '0.0982'
Set-TextValue 'E39' '  -1.22%  '
Set-TextValue 'E40' '  +0.68%  '
Set-TextValue 'D41' '15.26'
Set-TextValue 'E41' '  -4.94%  '
Set-TextValue 'D42' '2.111.94'
Set-TextValue 'E42' '  +2.51%  '
Set-TextValue 'D43' '3.73'
Set-TextValue 'E43' '  -1.72%  '
Set-TextValue 'D44' '0.0269'
Set-TextValue 'E44' '  +0.33%  '
Set-TextValue 'D45' '17.82'
Set-TextValue 'E45' '  -1.18%  '
Set-TextValue 'D46' '9.93'
Set-TextValue 'E46' '  -0.12%  '
Set-TextValue 'E47' '  -0.50%  '
Set-TextValue 'E48' '  +3.59%  '
Set-TextValue 'D49' '2.434.78'
Set-TextValue 'E49' '  +0.09%  '
Set-TextValue 'E50' '  +3.67%  '
Set-TextValue 'D51' '88.62'
Set-TextValue 'E51' '  -0.62%  '
